# Session 4 slide tweaks
#
# Helper: force PowerPoint's run-splitting engine to break a run at a given
# (0-based, .NET-string-index) offset by re-assigning a Characters() sub-range's
# Text to itself. PowerPoint (and this COM-interop emulation of it) only ever
# *creates* run boundaries where we touch a TextRange - touching a sub-range
# that starts mid-run forces a split there, while leaving runs that are fully
# outside the touched range untouched.
function Split-RunAt($textRange, [int]$offset0) {
    $len = $textRange.Length - $offset0
    if ($len -gt 0) {
        $c = $textRange.Characters($offset0 + 1, $len)
        $c.Text = $c.Text
    }
}

# Helper: force two (or more) adjacent runs that share identical formatting to
# collapse into a single run by re-assigning the Text of a Characters() range
# spanning all of them to itself.
function Merge-Runs($textRange, [int]$offset0, [int]$length) {
    $c = $textRange.Characters($offset0 + 1, $length)
    $c.Text = $c.Text
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1: title shape "Session 3: JavaScript - DOM and Events"
#   "3" -> "4"
#   split ": JavaScript - DOM and Events" into ": " and "JavaScript - DOM and Events"
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$titleShape = $s1.Shapes.Item(3)
$tr1 = $titleShape.TextFrame.TextRange
$txt1 = $tr1.Text

$idxSession = $txt1.IndexOf("Session ")
$sessionPrefixLen = "Session ".Length
$idxDigit = $idxSession + $sessionPrefixLen
$digitRange = $tr1.Characters($idxDigit + 1, 1)
$digitRange.Text = "4"

# re-read text (still same length, digit swap doesn't change offsets)
$txt1 = $tr1.Text
$idxJs = $txt1.IndexOf("JavaScript")
Split-RunAt $tr1 $idxJs

# ---------------------------------------------------------------------------
# Slide 10: "Get a node's / number of children:" and "Natural / to iterate..."
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$body10 = $s10.Shapes.Item(2)
$tr10 = $body10.TextFrame.TextRange
$txt10 = $tr10.Text

# NB: this merge's text contains a right single quotation mark (U+2019,
# "Get a node's"). The Characters().Text *getter* in this COM emulation
# normalizes that glyph down to a plain apostrophe, so a Text = Text
# round-trip here would silently corrupt the character. Assign the exact
# literal string instead so the curly apostrophe survives untouched.
$idxNumOf = $txt10.IndexOf("number of children")
$idxGetNode = $txt10.LastIndexOf("Get a node", $idxNumOf)
$idxListLen = $txt10.IndexOf("list.childNodes.length")
$mergeLen1 = $idxListLen - $idxGetNode
$getNodeRange = $tr10.Characters($idxGetNode + 1, $mergeLen1)
$getNodeRange.Text = "Get a node’s number of children: "

$txt10 = $tr10.Text
$idxIterate = $txt10.IndexOf("to iterate over child nodes")
$idxNatural = $txt10.LastIndexOf("Natural", $idxIterate)
$idxForLoops = $txt10.IndexOf("for loops")
$forLoopsLen = "for loops".Length
$endOfForLoops = $idxForLoops + $forLoopsLen
$mergeLen2 = $endOfForLoops - $idxNatural
Merge-Runs $tr10 $idxNatural $mergeLen2

# ---------------------------------------------------------------------------
# Slide 12: "Building / DOM nodes programmatically:"
# ---------------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$body12 = $s12.Shapes.Item(2)
$tr12 = $body12.TextFrame.TextRange
$txt12 = $tr12.Text

$idxBuilding = $txt12.IndexOf("Building DOM nodes programmatically")
$lenBuilding = "Building DOM nodes programmatically:".Length
Merge-Runs $tr12 $idxBuilding $lenBuilding

# ---------------------------------------------------------------------------
# Slide 13: " method / in its parent:"
# ---------------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$body13 = $s13.Shapes.Item(2)
$tr13 = $body13.TextFrame.TextRange
$txt13 = $tr13.Text

$idxMethod = $txt13.IndexOf(" method ")
$idxParent = $txt13.IndexOf("in its parent:")
$parentLen = "in its parent:".Length
$endOfParent = $idxParent + $parentLen
$mergeLen3 = $endOfParent - $idxMethod
Merge-Runs $tr13 $idxMethod $mergeLen3

# ---------------------------------------------------------------------------
# Slide 17: "Wikipedia / (Form (Document))"
# ---------------------------------------------------------------------------
$s17 = $p.Slides.Item(17)
$src17 = $s17.Shapes.Item(2)
$tr17 = $src17.TextFrame.TextRange
$txt17 = $tr17.Text

$idxWiki17 = $txt17.IndexOf("Wikipedia")
$idxFormDoc = $txt17.IndexOf("(Form (Document))")
$formDocLen = "(Form (Document))".Length
$idxEnd17 = $idxFormDoc + $formDocLen
$mergeLen4 = $idxEnd17 - $idxWiki17
Merge-Runs $tr17 $idxWiki17 $mergeLen4

# ---------------------------------------------------------------------------
# Slide 18: "Wikipedia / (Construction)"
# ---------------------------------------------------------------------------
$s18 = $p.Slides.Item(18)
$src18 = $s18.Shapes.Item(2)
$tr18 = $src18.TextFrame.TextRange
$txt18 = $tr18.Text

$idxWiki18 = $txt18.IndexOf("Wikipedia")
$idxConstr = $txt18.IndexOf("(Construction)")
$constrLen = "(Construction)".Length
$idxEnd18 = $idxConstr + $constrLen
$mergeLen5 = $idxEnd18 - $idxWiki18
Merge-Runs $tr18 $idxWiki18 $mergeLen5
